$d = $word.ActiveDocument

# 1. Update email address: "ctp8441" + "@rit.edu" runs -> single "cpittman343@gmail.com" run
$d.Content.Find.Execute("ctp8441@rit.edu", $false, $false, $false, $false, $false, $true, 1, $false, "cpittman343@gmail.com", 2) | Out-Null
